# Apply the "Testing all suitcases and changing xpaths" edit to the
# footer sheet: rename the first footer header label and rewrite the
# SelectValue xpaths to the shorter //*[@id="footer"]/... form.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("footer")

# Row 2: header label changes from footer_header_conocenos -> footer_header_producto
$ws.Range("A2").Value = "footer_header_producto"

# Rewrite the xpath (SelectValue) column to the //*[@id="footer"]/... form
$ws.Range("C2").Value  = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[1]'
$ws.Range("C3").Value  = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[5]/ul/li[1]'
$ws.Range("C4").Value  = '//*[@id="footer"]/div[1]/div[3]/div/div[2]/span'
$ws.Range("C5").Value  = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[2]/a'
$ws.Range("C6").Value  = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[2]/a'
$ws.Range("C7").Value  = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[3]/a'
$ws.Range("C8").Value  = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[3]/a'
$ws.Range("C9").Value  = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[4]/a'
$ws.Range("C10").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[4]/a'
$ws.Range("C11").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[5]/a'
$ws.Range("C12").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[2]/ul/li[5]/a'
$ws.Range("C13").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[5]/ul/li[2]/a'
$ws.Range("C14").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[5]/ul/li[2]/a'
$ws.Range("C15").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[5]/ul/li[3]/a'
$ws.Range("C16").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[5]/ul/li[3]/a'
$ws.Range("C17").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[5]/ul/li[4]/a'
$ws.Range("C18").Value = '//*[@id="footer"]/div[1]/div[2]/div[1]/div[5]/ul/li[4]/a'

# Move the active selection to C2 (matches the saved view state in the target file)
$ws.Range("C2").Select()
